$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: update text value (string "123456" -> "123"), preserving style/type ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "123"
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Updated numeric cells in row 2 ---
$ws.Range("C2").Value2 = [double]"2.714557972406106"
$ws.Range("D2").Value2 = [double]"53.61014775876133"
$ws.Range("E2").Value2 = [double]"31.85920189300484"
$ws.Range("F2").Value2 = [double]"4.044998102087909"
$ws.Range("G2").Value2 = [double]"0.8175396413048579"
$ws.Range("H2").Value2 = [double]"0.04905239148793471"
$ws.Range("I2").Value2 = [double]"6.880456211966506"
$ws.Range("J2").Value2 = [double]"0.006958268185923763"
$ws.Range("K2").Value2 = [double]"0.002918062911566768"
$ws.Range("M2").Value2 = [double]"0.0001172472693851281"
$ws.Range("N2").Value2 = [double]"0.01404498796819165"
$ws.Range("O2").Value2 = [double]"7.373085978330415e-06"
$ws.Range("P2").Value2 = [double]"2.85973636269423e-09"
$ws.Range("Q2").Value2 = [double]"1.305486610559374e-09"
$ws.Range("R2").Value2 = [double]"8.885468789419644e-09"
$ws.Range("S2").Value2 = [double]"3.007374882504019e-08"
$ws.Range("T2").Value2 = [double]"4.46840976340793e-08"
$ws.Range("U2").Value2 = [double]"1.717297498290396e-09"
$ws.Range("V2").Value2 = [double]"3.078999320135099e-11"
$ws.Range("W2").Value2 = [double]"7.040364470886498e-14"
$ws.Range("X2").Value2 = [double]"2.692172379612205e-14"
$ws.Range("Y2").Value2 = [double]"3.62001752650698e-15"
$ws.Range("Z2").Value2 = [double]"1.813018511007507e-12"
$ws.Range("AA2").Value2 = [double]"5.653869991586291e-13"
$ws.Range("AB2").Value2 = [double]"9.983457881886127e-14"
$ws.Range("AC2").Value2 = [double]"2.394024086302311e-13"
$ws.Range("AD2").Value2 = [double]"2.830152048668142e-14"
$ws.Range("AF2").Value2 = [double]"1.838838265108897e-16"
$ws.Range("AG2").Value2 = [double]"1.336674877886005e-17"
$ws.Range("AH2").Value2 = [double]"1.911714338071911e-20"
$ws.Range("AI2").Value2 = [double]"5.214238638919933e-19"
$ws.Range("AJ2").Value2 = [double]"1.067825299018216e-17"
$ws.Range("AL2").Value2 = [double]"6.017995180556357e-21"
$ws.Range("AM2").Value2 = [double]"1.197862228101574e-22"
$ws.Range("AN2").Value2 = [double]"2.783586674746736e-27"
$ws.Range("AR2").Value2 = [double]"0.02012458800113638"
$ws.Range("AS2").Value2 = [double]"2.202305264070173"
$ws.Range("AT2").Value2 = [double]"4.458863145824066"
$ws.Range("AU2").Value2 = [double]"2.008915703459953"
$ws.Range("AV2").Value2 = [double]"1.028977220943492"
$ws.Range("AW2").Value2 = [double]"0.1979250644895254"
$ws.Range("AX2").Value2 = [double]"88.18995186146206"
$ws.Range("AY2").Value2 = [double]"0.05567373943419723"
$ws.Range("AZ2").Value2 = [double]"0.05210050085371841"
$ws.Range("BB2").Value2 = [double]"0.004904410262486084"
$ws.Range("BC2").Value2 = [double]"1.779863368361432"
$ws.Range("BD2").Value2 = [double]"0.0002735723090076049"
$ws.Range("BE2").Value2 = [double]"1.587124630256042e-06"
$ws.Range("BF2").Value2 = [double]"7.508297550707308e-07"
$ws.Range("BG2").Value2 = [double]"1.039593591986309e-05"
$ws.Range("BH2").Value2 = [double]"4.058391001465674e-05"
$ws.Range("BI2").Value2 = [double]"6.492206874206049e-05"
$ws.Range("BJ2").Value2 = [double]"3.262685949879385e-06"
$ws.Range("BK2").Value2 = [double]"4.920432927120795e-08"
$ws.Range("BL2").Value2 = [double]"1.743123860794339e-10"
$ws.Range("BM2").Value2 = [double]"4.346027438755534e-11"
$ws.Range("BN2").Value2 = [double]"1.357736126419096e-11"
$ws.Range("BO2").Value2 = [double]"5.299477308747644e-09"
$ws.Range("BP2").Value2 = [double]"1.729649485529727e-09"
$ws.Range("BQ2").Value2 = [double]"6.079746922147797e-10"
$ws.Range("BR2").Value2 = [double]"6.931735888026281e-10"
$ws.Range("BS2").Value2 = [double]"2.056358526796096e-10"
$ws.Range("BU2").Value2 = [double]"1.916941504206999e-12"
$ws.Range("BV2").Value2 = [double]"9.628603074068967e-14"
$ws.Range("BW2").Value2 = [double]"2.389749052096902e-16"
$ws.Range("BX2").Value2 = [double]"5.113816315186424e-15"
$ws.Range("BY2").Value2 = [double]"1.506884660930347e-13"
$ws.Range("CA2").Value2 = [double]"2.143609086947379e-16"
$ws.Range("CB2").Value2 = [double]"5.215989601476102e-18"
$ws.Range("CC2").Value2 = [double]"3.783345908587044e-22"
$ws.Range("CG2").Value2 = [double]"0.0004208678838225773"
$ws.Range("CH2").Value2 = [double]"0.02725747827207207"
$ws.Range("CI2").Value2 = [double]"0.05168924795800026"
$ws.Range("CJ2").Value2 = [double]"0.0338216906862601"
$ws.Range("CK2").Value2 = [double]"0.03167619490045555"
$ws.Range("CL2").Value2 = [double]"0.02315386755688598"
$ws.Range("CM2").Value2 = [double]"87.43333469538887"
$ws.Range("CN2").Value2 = [double]"0.02757000931783608"
$ws.Range("CO2").Value2 = [double]"0.06714042325277385"
$ws.Range("CQ2").Value2 = [double]"0.02139042754268897"
$ws.Range("CR2").Value2 = [double]"12.27442077234601"
$ws.Range("CS2").Value2 = [double]"0.001415246326867387"
$ws.Range("CT2").Value2 = [double]"2.907958838606602e-05"
$ws.Range("CU2").Value2 = [double]"1.675278308497874e-05"
$ws.Range("CV2").Value2 = [double]"0.0004466863780926031"
$ws.Range("CW2").Value2 = [double]"0.002129134857489691"
$ws.Range("CX2").Value2 = [double]"0.003736633543057934"
$ws.Range("CY2").Value2 = [double]"0.0003440025253977484"
$ws.Range("CZ2").Value2 = [double]"2.690877273982389e-06"
$ws.Range("DA2").Value2 = [double]"4.031145025793708e-08"
$ws.Range("DB2").Value2 = [double]"1.083089214444096e-08"
$ws.Range("DC2").Value2 = [double]"6.995207608902696e-09"
$ws.Range("DD2").Value2 = [double]"1.994437677977014e-06"
$ws.Range("DE2").Value2 = [double]"7.46086003674642e-07"
$ws.Range("DF2").Value2 = [double]"5.957727141325702e-07"
$ws.Range("DG2").Value2 = [double]"3.288556059211307e-07"
$ws.Range("DH2").Value2 = [double]"3.651216312264359e-07"
$ws.Range("DJ2").Value2 = [double]"7.962141205368625e-09"
$ws.Range("DK2").Value2 = [double]"3.437798855245152e-10"
$ws.Range("DL2").Value2 = [double]"1.760992992909378e-12"
$ws.Range("DM2").Value2 = [double]"3.786675882453307e-11"
$ws.Range("DN2").Value2 = [double]"1.243992654359028e-09"
$ws.Range("DP2").Value2 = [double]"1.327562442819379e-11"
$ws.Range("DQ2").Value2 = [double]"6.388315408255418e-13"
$ws.Range("DR2").Value2 = [double]"6.236427960175826e-16"
$ws.Range("DY2").Value2 = [double]"9.263201066462393e-24"
$ws.Range("DZ2").Value2 = [double]"2.546577697377536e-16"
$ws.Range("EA2").Value2 = [double]"1.892381882233796e-09"
$ws.Range("EB2").Value2 = [double]"0.2256758726445471"
$ws.Range("EC2").Value2 = [double]"0.0001811135712232907"
$ws.Range("ED2").Value2 = [double]"0.003047838540434291"
$ws.Range("EF2").Value2 = [double]"0.009122709887292082"
$ws.Range("EG2").Value2 = [double]"45.76374715986674"
$ws.Range("EH2").Value2 = [double]"0.008010068768777851"
$ws.Range("EI2").Value2 = [double]"0.009670091920880119"
$ws.Range("EJ2").Value2 = [double]"0.009671803910525661"
$ws.Range("EK2").Value2 = [double]"1.145711807743856"
$ws.Range("EL2").Value2 = [double]"8.947271040397032"
$ws.Range("EM2").Value2 = [double]"19.54510178117068"
$ws.Range("EN2").Value2 = [double]"8.006474343521461"
$ws.Range("EO2").Value2 = [double]"0.01403993070888112"
$ws.Range("EP2").Value2 = [double]"0.006755792581057265"
$ws.Range("EQ2").Value2 = [double]"0.002653956544538992"
$ws.Range("ER2").Value2 = [double]"0.007881762802523697"
$ws.Range("ES2").Value2 = [double]"1.104572754381669"
$ws.Range("ET2").Value2 = [double]"0.578746568468153"
$ws.Range("EU2").Value2 = [double]"2.932016215395487"
$ws.Range("EV2").Value2 = [double]"0.3254042076530504"
$ws.Range("EW2").Value2 = [double]"7.488801002821739"
$ws.Range("EY2").Value2 = [double]"1.195776195305765"
$ws.Range("EZ2").Value2 = [double]"0.03940881758374768"
$ws.Range("FA2").Value2 = [double]"0.001125966279057494"
$ws.Range("FB2").Value2 = [double]"0.02589722420556916"
$ws.Range("FC2").Value2 = [double]"1.01224357672379"
$ws.Range("FE2").Value2 = [double]"1.164249199079529"
$ws.Range("FF2").Value2 = [double]"0.2983810371040776"
$ws.Range("FG2").Value2 = [double]"0.1283601585255211"
